$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028494635305535
$ws.Cells.Item(2, 4).Value = 1.029744293220024
$ws.Cells.Item(2, 5).Value = 1.037409224207617
$ws.Cells.Item(2, 6).Value = 1.045952842655858
$ws.Cells.Item(2, 9).Value = 1.026982128717641
$ws.Cells.Item(2, 10).Value = 1.033646640152938
$ws.Cells.Item(2, 11).Value = 1.032557315927483
$ws.Cells.Item(2, 12).Value = 1.040200191076355
$ws.Cells.Item(2, 13).Value = 1.048719638666159
$ws.Cells.Item(2, 14).Value = 1.035114536942882

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029756007235823
$ws.Cells.Item(3, 4).Value = 1.030853427006506
$ws.Cells.Item(3, 5).Value = 1.038528627343567
$ws.Cells.Item(3, 6).Value = 1.047157444064425
$ws.Cells.Item(3, 9).Value = 1.026958505582768
$ws.Cells.Item(3, 10).Value = 1.034547047154769
$ws.Cells.Item(3, 11).Value = 1.033473849825907
$ws.Cells.Item(3, 12).Value = 1.041128576365803
$ws.Cells.Item(3, 13).Value = 1.049734762505478
$ws.Cells.Item(3, 14).Value = 1.036016222625933

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030572211910802
$ws.Cells.Item(4, 4).Value = 1.031571387647306
$ws.Cells.Item(4, 5).Value = 1.039253043647128
$ws.Cells.Item(4, 6).Value = 1.047936460196186
$ws.Cells.Item(4, 9).Value = 1.026940799041256
$ws.Cells.Item(4, 10).Value = 1.035129224846839
$ws.Cells.Item(4, 11).Value = 1.03406660368445
$ws.Cells.Item(4, 12).Value = 1.041728809471719
$ws.Cells.Item(4, 13).Value = 1.050390608920692
$ws.Cells.Item(4, 14).Value = 1.036599227077103

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030915349706467
$ws.Cells.Item(5, 4).Value = 1.031873285871108
$ws.Cells.Item(5, 5).Value = 1.039557610383361
$ws.Cells.Item(5, 6).Value = 1.048263853896648
$ws.Cells.Item(5, 9).Value = 1.026932774768032
$ws.Cells.Item(5, 10).Value = 1.0353738668146
$ws.Cells.Item(5, 11).Value = 1.034315725586305
$ws.Cells.Item(5, 12).Value = 1.041981030161588
$ws.Cells.Item(5, 13).Value = 1.050666086534226
$ws.Cells.Item(5, 14).Value = 1.036844216464519

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030972964433163
$ws.Cells.Item(6, 4).Value = 1.031923979910697
$ws.Cells.Item(6, 5).Value = 1.039608749808981
$ws.Cells.Item(6, 6).Value = 1.048318818630509
$ws.Cells.Item(6, 9).Value = 1.026931393393354
$ws.Cells.Item(6, 10).Value = 1.035414937123782
$ws.Cells.Item(6, 11).Value = 1.034357550072004
$ws.Cells.Item(6, 12).Value = 1.042023372275737
$ws.Cells.Item(6, 13).Value = 1.050712326373008
$ws.Cells.Item(6, 14).Value = 1.036885345098252

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03057679691093
$ws.Cells.Item(7, 4).Value = 1.031575421359072
$ws.Cells.Item(7, 5).Value = 1.039257113194444
$ws.Cells.Item(7, 6).Value = 1.047940835256951
$ws.Cells.Item(7, 9).Value = 1.026940694102428
$ws.Cells.Item(7, 10).Value = 1.035132494177989
$ws.Cells.Item(7, 11).Value = 1.034069932744272
$ws.Cells.Item(7, 12).Value = 1.041732180116628
$ws.Cells.Item(7, 13).Value = 1.050394290809066
$ws.Cells.Item(7, 14).Value = 1.036602501051078

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028920920619624
$ws.Cells.Item(8, 4).Value = 1.030119073608446
$ws.Cells.Item(8, 5).Value = 1.037787514286764
$ws.Cells.Item(8, 6).Value = 1.046360035381099
$ws.Cells.Item(8, 9).Value = 1.026974645872505
$ws.Cells.Item(8, 10).Value = 1.033951029955483
$ws.Cells.Item(8, 11).Value = 1.032867126674882
$ws.Cells.Item(8, 12).Value = 1.040514046184984
$ws.Cells.Item(8, 13).Value = 1.049062913401597
$ws.Cells.Item(8, 14).Value = 1.035419359013861

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026003048490246
$ws.Cells.Item(9, 4).Value = 1.027554866008168
$ws.Cells.Item(9, 5).Value = 1.035198512425166
$ws.Cells.Item(9, 6).Value = 1.043571035760224
$ws.Cells.Item(9, 9).Value = 1.027015971983489
$ws.Cells.Item(9, 10).Value = 1.031865670099288
$ws.Cells.Item(9, 11).Value = 1.030745241268708
$ws.Cells.Item(9, 12).Value = 1.038363706976997
$ws.Cells.Item(9, 13).Value = 1.046709112046345
$ws.Cells.Item(9, 14).Value = 1.033331037707476

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024057648894561
$ws.Cells.Item(10, 4).Value = 1.025846698304632
$ws.Cells.Item(10, 5).Value = 1.033472851991469
$ws.Cells.Item(10, 6).Value = 1.041709328632811
$ws.Cells.Item(10, 9).Value = 1.027031135448596
$ws.Cells.Item(10, 10).Value = 1.030473019687299
$ws.Cells.Item(10, 11).Value = 1.029328970030086
$ws.Cells.Item(10, 12).Value = 1.036927504133633
$ws.Cells.Item(10, 13).Value = 1.045134651707326
$ws.Cells.Item(10, 14).Value = 1.031936409572162

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023215200068402
$ws.Cells.Item(11, 4).Value = 1.02510733282491
$ws.Cells.Item(11, 5).Value = 1.032725685099837
$ws.Cells.Item(11, 6).Value = 1.040902609248143
$ws.Cells.Item(11, 9).Value = 1.027034773977178
$ws.Cells.Item(11, 10).Value = 1.029869396833032
$ws.Cells.Item(11, 11).Value = 1.028715292661357
$ws.Cells.Item(11, 12).Value = 1.036304969680405
$ws.Cells.Item(11, 13).Value = 1.044451632662756
$ws.Cells.Item(11, 14).Value = 1.031331929504206

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022902262387607
$ws.Cells.Item(12, 4).Value = 1.024832740238272
$ws.Cells.Item(12, 5).Value = 1.032448160751356
$ws.Cells.Item(12, 6).Value = 1.040602867859209
$ws.Cells.Item(12, 9).Value = 1.027035686501887
$ws.Cells.Item(12, 10).Value = 1.029645093230966
$ws.Cells.Item(12, 11).Value = 1.028487280471654
$ws.Cells.Item(12, 12).Value = 1.036073633677478
$ws.Cells.Item(12, 13).Value = 1.044197737330244
$ws.Cells.Item(12, 14).Value = 1.031107307365296

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022969389294246
$ws.Cells.Item(13, 4).Value = 1.024891639467976
$ws.Cells.Item(13, 5).Value = 1.032507690389141
$ws.Cells.Item(13, 6).Value = 1.040667167497579
$ws.Cells.Item(13, 9).Value = 1.027035510614369
$ws.Cells.Item(13, 10).Value = 1.029693211268479
$ws.Cells.Item(13, 11).Value = 1.028536192831676
$ws.Cells.Item(13, 12).Value = 1.036123260538399
$ws.Cells.Item(13, 13).Value = 1.044252207433541
$ws.Cells.Item(13, 14).Value = 1.031155493735943

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023189332858654
$ws.Cells.Item(14, 4).Value = 1.025084634096511
$ws.Cells.Item(14, 5).Value = 1.032702744723347
$ws.Cells.Item(14, 6).Value = 1.040877834374252
$ws.Cells.Item(14, 9).Value = 1.027034858357287
$ws.Cells.Item(14, 10).Value = 1.029850857697288
$ws.Cells.Item(14, 11).Value = 1.028696446439765
$ws.Cells.Item(14, 12).Value = 1.036285849402063
$ws.Cells.Item(14, 13).Value = 1.044430649512182
$ws.Cells.Item(14, 14).Value = 1.031313364040763

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023324845314153
$ws.Cells.Item(15, 4).Value = 1.025203549793659
$ws.Cells.Item(15, 5).Value = 1.032822924975505
$ws.Cells.Item(15, 6).Value = 1.041007621231541
$ws.Cells.Item(15, 9).Value = 1.027034398334308
$ws.Cells.Item(15, 10).Value = 1.029947976731006
$ws.Cells.Item(15, 11).Value = 1.028795175303468
$ws.Cells.Item(15, 12).Value = 1.036386012606237
$ws.Cells.Item(15, 13).Value = 1.044540568133455
$ws.Cells.Item(15, 14).Value = 1.031410620994648

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024113557553223
$ws.Cells.Item(16, 4).Value = 1.025895773356799
$ws.Cells.Item(16, 5).Value = 1.033522440026956
$ws.Cells.Item(16, 6).Value = 1.041762855421732
$ws.Cells.Item(16, 9).Value = 1.027030832365333
$ws.Cells.Item(16, 10).Value = 1.030513067449685
$ws.Cells.Item(16, 11).Value = 1.02936968873603
$ws.Cells.Item(16, 12).Value = 1.036968805893698
$ws.Cells.Item(16, 13).Value = 1.045179954624613
$ws.Cells.Item(16, 14).Value = 1.031976514206965

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.0246082729397
$ws.Cells.Item(17, 4).Value = 1.026330061280549
$ws.Cells.Item(17, 5).Value = 1.033961241038419
$ws.Cells.Item(17, 6).Value = 1.042236435375054
$ws.Cells.Item(17, 9).Value = 1.027027812371069
$ws.Cells.Item(17, 10).Value = 1.030867373189134
$ws.Cells.Item(17, 11).Value = 1.029729951567668
$ws.Cells.Item(17, 12).Value = 1.037334201901263
$ws.Cells.Item(17, 13).Value = 1.045580684701939
$ws.Cells.Item(17, 14).Value = 1.032331323101203

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024896824734028
$ws.Cells.Item(18, 4).Value = 1.026583401447354
$ws.Cells.Item(18, 5).Value = 1.034217191707769
$ws.Cells.Item(18, 6).Value = 1.042512609932154
$ws.Cells.Item(18, 9).Value = 1.027025768358815
$ws.Cells.Item(18, 10).Value = 1.031073976336591
$ws.Cells.Item(18, 11).Value = 1.029940046119396
$ws.Cells.Item(18, 12).Value = 1.037547268677735
$ws.Cells.Item(18, 13).Value = 1.045814301583703
$ws.Cells.Item(18, 14).Value = 1.032538219648827

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024995212286415
$ws.Cells.Item(19, 4).Value = 1.02666978859204
$ws.Cells.Item(19, 5).Value = 1.034304465247008
$ws.Cells.Item(19, 6).Value = 1.042606768762347
$ws.Cells.Item(19, 9).Value = 1.027025023461371
$ws.Cells.Item(19, 10).Value = 1.031144412984824
$ws.Cells.Item(19, 11).Value = 1.030011676083092
$ws.Cells.Item(19, 12).Value = 1.037619908379397
$ws.Cells.Item(19, 13).Value = 1.04589393821417
$ws.Cells.Item(19, 14).Value = 1.03260875632518

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024555195440593
$ws.Cells.Item(20, 4).Value = 1.026283463479275
$ws.Cells.Item(20, 5).Value = 1.033914161290414
$ws.Cells.Item(20, 6).Value = 1.042185630596673
$ws.Cells.Item(20, 9).Value = 1.027028165602469
$ws.Cells.Item(20, 10).Value = 1.030829365486001
$ws.Cells.Item(20, 11).Value = 1.029691302993788
$ws.Cells.Item(20, 12).Value = 1.037295004847296
$ws.Cells.Item(20, 13).Value = 1.045537702806141
$ws.Cells.Item(20, 14).Value = 1.032293261422772

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.02312456538987
$ws.Cells.Item(21, 4).Value = 1.025027800880042
$ws.Cells.Item(21, 5).Value = 1.032645305911632
$ws.Cells.Item(21, 6).Value = 1.04081580072383
$ws.Cells.Item(21, 9).Value = 1.027035062542446
$ws.Cells.Item(21, 10).Value = 1.029804437284025
$ws.Cells.Item(21, 11).Value = 1.028649257556139
$ws.Cells.Item(21, 12).Value = 1.036237973778487
$ws.Cells.Item(21, 13).Value = 1.044378108068063
$ws.Cells.Item(21, 14).Value = 1.031266877705188

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.02222498510771
$ws.Cells.Item(22, 4).Value = 1.024238550412036
$ws.Cells.Item(22, 5).Value = 1.031847564140198
$ws.Cells.Item(22, 6).Value = 1.039954013453757
$ws.Cells.Item(22, 9).Value = 1.027036859474112
$ws.Cells.Item(22, 10).Value = 1.029159495780375
$ws.Cells.Item(22, 11).Value = 1.027993704835483
$ws.Cells.Item(22, 12).Value = 1.035572802944274
$ws.Cells.Item(22, 13).Value = 1.043647915390114
$ws.Cells.Item(22, 14).Value = 1.030621020310627

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022701878816376
$ws.Cells.Item(23, 4).Value = 1.024656925384118
$ws.Cells.Item(23, 5).Value = 1.032270459095257
$ws.Cells.Item(23, 6).Value = 1.04041091308833
$ws.Cells.Item(23, 9).Value = 1.027036147328315
$ws.Cells.Item(23, 10).Value = 1.029501442119508
$ws.Cells.Item(23, 11).Value = 1.028341262046446
$ws.Cells.Item(23, 12).Value = 1.03592547743583
$ws.Cells.Item(23, 13).Value = 1.044035109926381
$ws.Cells.Item(23, 14).Value = 1.030963452252783

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024579178900318
$ws.Cells.Item(24, 4).Value = 1.026304518935384
$ws.Cells.Item(24, 5).Value = 1.033935434585204
$ws.Cells.Item(24, 6).Value = 1.042208587264188
$ws.Cells.Item(24, 9).Value = 1.027028006865416
$ws.Cells.Item(24, 10).Value = 1.030846539708394
$ws.Cells.Item(24, 11).Value = 1.029708766746392
$ws.Cells.Item(24, 12).Value = 1.037312716502033
$ws.Cells.Item(24, 13).Value = 1.045557124852365
$ws.Cells.Item(24, 14).Value = 1.032310460034531

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026757404279123
$ws.Cells.Item(25, 4).Value = 1.028217539173205
$ws.Cells.Item(25, 5).Value = 1.035867765831303
$ws.Cells.Item(25, 6).Value = 1.04429247234877
$ws.Cells.Item(25, 9).Value = 1.027007475287657
$ws.Cells.Item(25, 10).Value = 1.032405204722763
$ws.Cells.Item(25, 11).Value = 1.031294090211725
$ws.Cells.Item(25, 12).Value = 1.038920082161394
$ws.Cells.Item(25, 13).Value = 1.047318549225851
$ws.Cells.Item(25, 14).Value = 1.033871338532001
